$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.605999999999992
$ws.Range("A8").Value = -21.17270000000001
$ws.Range("A10").Value = -20.47549999999998
$ws.Range("A12").Value = -22.44960000000003
$ws.Range("C13").Value = -13.4038
$ws.Range("A18").Value = -22.35560000000003
$ws.Range("D20").Value = -8.199700000000004
